# Update the cryptocurrency Price/Volume(1h) snapshot (and, for rows 48-49,
# swap which coin - Elrond vs EnergySwap - occupies which row) to match the
# latest GitHub Actions data pull.
#
# Column D ("Price") cells store plain text such as "1.720" or "0.2590"
# using dot-grouped formatting that is not valid numeric/date data. These
# cells have no explicit number format (General), so assigning the bare
# string to .Value would let Excel auto-detect it as a number and silently
# drop meaningful trailing zeros (e.g. "1.720" -> 1.72). Prefixing the
# literal with a leading apostrophe forces Excel to store it as text,
# exactly as a user typing these values manually would do to keep them
# intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.384.52'
$ws.Range("D3").Value = '''1.722.65'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").Value = '''0.9993'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''242.73'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '''0.4878'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '''0.2590'
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = '''0.06190'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''1.739.53'
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").Value = '''0.06978'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = '''4.531'
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '''0.5969'
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("D15").Value = '''77.11'
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = '''0.9998'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '''26.381.93'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '''0.9995'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '''0.000007185'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = '''11.33'
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("D21").Value = '''1.947.59'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  -1.80%  '
$ws.Range("D23").Value = '''8.488'
$ws.Range("E23").Value = '  -3.52%  '
$ws.Range("D24").Value = '''5.104'
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("D25").Value = '''137.97'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '''15.23'
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("D27").Value = '''1.398'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").Value = '''106.72'
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("D29").Value = '''1.720'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("D30").Value = '''3.918'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("D31").Value = '''0.07999'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '''3.665'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").Value = '''0.04501'
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("D34").Value = '''2.605'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Value = '''0.9958'
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").Value = '''0.6241'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '''0.9313'
$ws.Range("E37").Value = '  +3.71%  '
$ws.Range("D38").Value = '''1.959'
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("D39").Value = '''2.389'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''0.9990'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").Value = '''0.01472'
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("D42").Value = '''100.06'
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").Value = '''5.462'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = '''0.3836'
$ws.Range("E44").Value = '  -1.80%  '
$ws.Range("D45").Value = '''6.844'
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("D46").Value = '''0.1162'
$ws.Range("E46").Value = '  -2.06%  '
$ws.Range("D47").Value = '''0.05364'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''7.715'
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.06'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("D51").Value = '''50.84'
$ws.Range("E51").Value = '  -1.31%  '
